# Update the cryptocurrency price/volume table with the latest scraped
# values (GitHub Actions crypto-list refresh).
#
# Columns D (Price) and E (Volume 1h) on Sheet1 hold text strings (not
# numbers) -- several prices use '.' as a thousands separator (e.g.
# "43.709.64") and must stay literal text. For prices that happen to look
# like a plain decimal number (e.g. "239.08"), a leading apostrophe forces
# Excel to keep the cell as text instead of silently converting it to a
# floating point number (which would also lose meaningful trailing zeros,
# e.g. "10.30" -> 10.3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.709.64'
$ws.Range("E2").Value = '  -0.46%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '''239.08'
$ws.Range("E5").Value = '  -1.13%  '
$ws.Range("D6").Value = '''0.663'
$ws.Range("E6").Value = '  -4.07%  '
$ws.Range("D7").Value = '''72.05'
$ws.Range("E7").Value = '  -6.55%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '''0.593'
$ws.Range("E9").Value = '  -6.92%  '
$ws.Range("D10").Value = '''0.0984'
$ws.Range("E10").Value = '  -4.31%  '
$ws.Range("D11").Value = '''58.25'
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("E12").Value = '  -3.80%  '
$ws.Range("D13").Value = '''0.107'
$ws.Range("E13").Value = '  -0.78%  '
$ws.Range("D14").Value = '''7.12'
$ws.Range("E14").Value = '  -6.45%  '
$ws.Range("D15").Value = '2.683.63'
$ws.Range("E15").Value = '  -1.12%  '
$ws.Range("D16").Value = '''16.03'
$ws.Range("E16").Value = '  -5.31%  '
$ws.Range("D17").Value = '''0.895'
$ws.Range("E17").Value = '  -3.13%  '
$ws.Range("D18").Value = '2.334.46'
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("D19").Value = '43.628.05'
$ws.Range("E19").Value = '  -0.54%  '
$ws.Range("E20").Value = '  -3.34%  '
$ws.Range("D21").Value = '''77.86'
$ws.Range("D22").Value = '''6.57'
$ws.Range("E22").Value = '  -1.76%  '
$ws.Range("D23").Value = '''250.58'
$ws.Range("E23").Value = '  -2.53%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").Value = '''1.89'
$ws.Range("E25").Value = '  +6.27%  '
$ws.Range("E26").Value = '  +2.61%  '
$ws.Range("E27").Value = '  -2.22%  '
$ws.Range("D28").Value = '''10.30'
$ws.Range("E28").Value = '  -7.52%  '
$ws.Range("D29").Value = '''2.27'
$ws.Range("E29").Value = '  -1.09%  '
$ws.Range("D30").Value = '''175.72'
$ws.Range("E30").Value = '  +0.67%  '
$ws.Range("D31").Value = '''22.12'
$ws.Range("E31").Value = '  -4.41%  '
$ws.Range("E32").Value = '  -2.66%  '
$ws.Range("D33").Value = '''0.134'
$ws.Range("E33").Value = '  -1.72%  '
$ws.Range("D34").Value = '''0.0732'
$ws.Range("E34").Value = '  -3.27%  '
$ws.Range("D35").Value = '''5.05'
$ws.Range("E35").Value = '  -5.30%  '
$ws.Range("D36").Value = '''5.33'
$ws.Range("E36").Value = '  -1.55%  '
$ws.Range("E37").Value = '  -2.69%  '
$ws.Range("E38").Value = '  -1.65%  '
$ws.Range("D39").Value = '''2.36'
$ws.Range("E39").Value = '  -3.71%  '
$ws.Range("D40").Value = '''5.71'
$ws.Range("E40").Value = '  +26.52%  '
$ws.Range("D41").Value = '''0.0270'
$ws.Range("E41").Value = '  -2.94%  '
$ws.Range("D42").Value = '''65.25'
$ws.Range("E42").Value = '  +17.76%  '
$ws.Range("D43").Value = '''9.19'
$ws.Range("E43").Value = '  +2.19%  '
$ws.Range("E44").Value = '  +4.15%  '
$ws.Range("D45").Value = '''18.62'
$ws.Range("E45").Value = '  -3.33%  '
$ws.Range("D46").Value = '''0.196'
$ws.Range("E46").Value = '  -3.85%  '
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("D48").Value = '''2.98'
$ws.Range("E48").Value = '  +5.56%  '
$ws.Range("D49").Value = '''1.22'
$ws.Range("E49").Value = '  -3.87%  '
$ws.Range("D50").Value = '''2.41'
$ws.Range("E50").Value = '  -4.42%  '
$ws.Range("D51").Value = '''97.69'
$ws.Range("E51").Value = '  -4.29%  '
